$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Berekening oversterfte")

# --- Update individual "Waargenomen" (G) and "Verwacht" (H) values for existing rows ---
# Row 3
$ws.Range("G3").Value = 3217
# Row 4
$ws.Range("G4").Value = 3615
# Row 5
$ws.Range("G5").Value = 4459
# Row 7
$ws.Range("H7").Value = 2909
# Row 8
$ws.Range("H8").Value = 3010
# Row 9
$ws.Range("G9").Value = 3906
# Row 11
$ws.Range("G11").Value = 2984
$ws.Range("H11").Value = 2933
# Row 12
$ws.Range("H12").Value = 3050
# Row 13
$ws.Range("G13").Value = 2770
# Row 14
$ws.Range("G14").Value = 2725
# Row 17
$ws.Range("G17").Value = 2692
# Row 19
$ws.Range("G19").Value = 2636
# Row 20
$ws.Range("G20").Value = 2614
$ws.Range("H20").Value = 2856
# Row 21
$ws.Range("G21").Value = 2526
# Row 22
$ws.Range("G22").Value = 2670
# Row 23
$ws.Range("G23").Value = 2657
# Row 24
$ws.Range("G24").Value = 2634
# Row 25
$ws.Range("G25").Value = 3202
# Row 26
$ws.Range("G26").Value = 2836
# Row 27
$ws.Range("G27").Value = 2715
# Row 28
$ws.Range("G28").Value = 2661

# --- Insert a new row 29 (week 37) above the old totals row, pushing it from row 30 to row 31 ---
$ws.Rows.Item(29).Insert()

$ws.Range("F29").Value = 37
$ws.Range("G29").Value = 2704
$ws.Range("H29").Value = 2844
$ws.Range("I29").Formula = "=G29-H29"

# --- Update the sheet view: scroll back to top-left and select I13 ---
$ws.Range("A1").Select()
$ws.Range("I13").Select()
